$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wrap the three hint/reasoning lines in green color tags (matches the in-game
# "thought" styling) and add the leading space that the author introduced.
$ws.Range("B2").Value = " <color=#00CC00>(When you pointed out the killer just now, one person reacted rather unnaturally.)</color>"
$ws.Range("B3").Value = " <color=#00CC00>(Think back to this person’s testimony from last night.)</color>"
$ws.Range("B4").Value = " <color=#00CC00>(Try to spot the contradiction in this person’s statement.)</color>"

# The longer text now wraps onto a second line in column B (wrapText style),
# so rows 3 and 4 need to grow to match row 2's two-line height.
$ws.Rows.Item(3).RowHeight = 34
$ws.Rows.Item(4).RowHeight = 34

# Move the saved cursor/selection to B12, as in the authored file.
$ws.Range("B12").Select()
